$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 151.5090713333333
$ws.Range("H2").Value = 454.527214
$ws.Range("I2").Value = 0.7011959105080696
$ws.Range("J2").Value = 0.7120250504896967
$ws.Range("M2").Value = 162.7225033333333
$ws.Range("N2").Value = 488.16751
$ws.Range("O2").Value = 0.5231437953541009
$ws.Range("P2").Value = 0.5247717033381212
$ws.Range("Q2").Value = 24653.93536506857
$ws.Range("R2").Value = 221885.4182856171
$ws.Range("S2").Value = 0.3668262899099661
$ws.Range("T2").Value = 0.3736505985648899
$ws.Range("G3").Value = 151.5090713333333
$ws.Range("H3").Value = 454.527214
$ws.Range("I3").Value = 0.7011959105080696
$ws.Range("J3").Value = 0.7120250504896967
$ws.Range("O3").Value = 0.0009322191998643353
$ws.Range("P3").Value = 0.0009351200601857102
$ws.Range("Q3").Value = 43.93222686312244
$ws.Range("R3").Value = 395.390041768102
$ws.Range("S3").Value = 0.0006536682906419768
$ws.Range("T3").Value = 0.0006658289080676584
$ws.Range("G4").Value = 151.5090713333333
$ws.Range("H4").Value = 454.527214
$ws.Range("I4").Value = 0.7011959105080696
$ws.Range("J4").Value = 0.7120250504896967
$ws.Range("M4").Value = 61.580654
$ws.Range("N4").Value = 184.741962
$ws.Range("O4").Value = 0.1979783766474813
$ws.Range("P4").Value = 0.1985944416431287
$ws.Range("Q4").Value = 9330.027699639319
$ws.Range("R4").Value = 83970.24929675387
$ws.Range("S4").Value = 0.1388216280742402
$ws.Range("T4").Value = 0.1414042173379218
$ws.Range("G5").Value = 151.5090713333333
$ws.Range("H5").Value = 454.527214
$ws.Range("I5").Value = 0.7011959105080696
$ws.Range("J5").Value = 0.7120250504896967
$ws.Range("M5").Value = 2.8947245
$ws.Range("N5").Value = 5.789449
$ws.Range("O5").Value = 0.009306378223129816
$ws.Range("P5").Value = 0.00622355841157717
$ws.Range("Q5").Value = 438.5770207608477
$ws.Range("R5").Value = 2631.462124565086
$ws.Range("S5").Value = 0.006525594351699983
$ws.Range("T5").Value = 0.004431329492228811
$ws.Range("G6").Value = 151.5090713333333
$ws.Range("H6").Value = 454.527214
$ws.Range("I6").Value = 0.7011959105080696
$ws.Range("J6").Value = 0.7120250504896967
$ws.Range("M6").Value = 83.559527
$ws.Range("N6").Value = 250.678581
$ws.Range("O6").Value = 0.2686392305754237
$ws.Range("P6").Value = 0.2694751765469873
$ws.Range("Q6").Value = 12660.02633682259
$ws.Range("R6").Value = 113940.2370314033
$ws.Range("S6").Value = 0.1883687298815215
$ws.Range("T6").Value = 0.1918730761865885
$ws.Range("I7").Value = 0.2398807347813626
$ws.Range("J7").Value = 0.2435854084922527
$ws.Range("M7").Value = 162.7225033333333
$ws.Range("N7").Value = 488.16751
$ws.Range("O7").Value = 0.5231437953541009
$ws.Range("P7").Value = 0.5247717033381212
$ws.Range("Q7").Value = 8434.168029217006
$ws.Range("R7").Value = 75907.51226295305
$ws.Range("S7").Value = 0.1254921180258525
$ws.Range("T7").Value = 0.1278267297227915
$ws.Range("I8").Value = 0.2398807347813626
$ws.Range("J8").Value = 0.2435854084922527
$ws.Range("O8").Value = 0.0009322191998643353
$ws.Range("P8").Value = 0.0009351200601857102
$ws.Range("S8").Value = 0.0002236214266407507
$ws.Range("T8").Value = 0.0002277816018496362
$ws.Range("I9").Value = 0.2398807347813626
$ws.Range("J9").Value = 0.2435854084922527
$ws.Range("M9").Value = 61.580654
$ws.Range("N9").Value = 184.741962
$ws.Range("O9").Value = 0.1979783766474813
$ws.Range("P9").Value = 0.1985944416431287
$ws.Range("Q9").Value = 3191.823949027708
$ws.Range("R9").Value = 28726.41554124937
$ws.Range("S9").Value = 0.04749119846101917
$ws.Range("T9").Value = 0.04837470819193235
$ws.Range("I10").Value = 0.2398807347813626
$ws.Range("J10").Value = 0.2435854084922527
$ws.Range("M10").Value = 2.8947245
$ws.Range("N10").Value = 5.789449
$ws.Range("O10").Value = 0.009306378223129816
$ws.Range("P10").Value = 0.00622355841157717
$ws.Range("Q10").Value = 150.038208183649
$ws.Range("R10").Value = 900.229249101894
$ws.Range("S10").Value = 0.002232420846317652
$ws.Range("T10").Value = 0.001515968017959421
$ws.Range("I11").Value = 0.2398807347813626
$ws.Range("J11").Value = 0.2435854084922527
$ws.Range("M11").Value = 83.559527
$ws.Range("N11").Value = 250.678581
$ws.Range("O11").Value = 0.2686392305754237
$ws.Range("P11").Value = 0.2694751765469873
$ws.Range("Q11").Value = 4331.024146772254
$ws.Range("R11").Value = 38979.21732095028
$ws.Range("S11").Value = 0.06444137602153251
$ws.Range("T11").Value = 0.06564022095771982
$ws.Range("G12").Value = 1.349157666666667
$ws.Range("H12").Value = 4.047473
$ws.Range("I12").Value = 0.006244007901123888
$ws.Range("J12").Value = 0.006340439204550433
$ws.Range("M12").Value = 162.7225033333333
$ws.Range("N12").Value = 488.16751
$ws.Range("O12").Value = 0.5231437953541009
$ws.Range("P12").Value = 0.5247717033381212
$ws.Range("Q12").Value = 219.5383129113589
$ws.Range("R12").Value = 1975.84481620223
$ws.Range("S12").Value = 0.003266513991614944
$ws.Range("T12").Value = 0.003327283081283733
$ws.Range("G13").Value = 1.349157666666667
$ws.Range("H13").Value = 4.047473
$ws.Range("I13").Value = 0.006244007901123888
$ws.Range("J13").Value = 0.006340439204550433
$ws.Range("O13").Value = 0.0009322191998643353
$ws.Range("P13").Value = 0.0009351200601857102
$ws.Range("Q13").Value = 0.3912076033765556
$ws.Range("R13").Value = 3.520868430389
$ws.Range("S13").Value = 0.000005820784049532299
$ws.Range("T13").Value = 0.000005929071890563037
$ws.Range("G14").Value = 1.349157666666667
$ws.Range("H14").Value = 4.047473
$ws.Range("I14").Value = 0.006244007901123888
$ws.Range("J14").Value = 0.006340439204550433
$ws.Range("M14").Value = 61.580654
$ws.Range("N14").Value = 184.741962
$ws.Range("O14").Value = 0.1979783766474813
$ws.Range("P14").Value = 0.1985944416431287
$ws.Range("Q14").Value = 83.08201146244734
$ws.Range("R14").Value = 747.7381031620261
$ws.Range("S14").Value = 0.001236178548038554
$ws.Range("T14").Value = 0.001259175983599896
$ws.Range("G15").Value = 1.349157666666667
$ws.Range("H15").Value = 4.047473
$ws.Range("I15").Value = 0.006244007901123888
$ws.Range("J15").Value = 0.006340439204550433
$ws.Range("M15").Value = 2.8947245
$ws.Range("N15").Value = 5.789449
$ws.Range("O15").Value = 0.009306378223129816
$ws.Range("P15").Value = 0.00622355841157717
$ws.Range("Q15").Value = 3.905439752062834
$ws.Range("R15").Value = 23.432638512377
$ws.Range("S15").Value = 0.00005810909915606985
$ws.Range("T15").Value = 0.00003946009374457351
$ws.Range("G16").Value = 1.349157666666667
$ws.Range("H16").Value = 4.047473
$ws.Range("I16").Value = 0.006244007901123888
$ws.Range("J16").Value = 0.006340439204550433
$ws.Range("M16").Value = 83.559527
$ws.Range("N16").Value = 250.678581
$ws.Range("O16").Value = 0.2686392305754237
$ws.Range("P16").Value = 0.2694751765469873
$ws.Range("Q16").Value = 112.7349764750903
$ws.Range("R16").Value = 1014.614788275813
$ws.Range("S16").Value = 0.001677385478264787
$ws.Range("T16").Value = 0.001708590974031667
$ws.Range("G17").Value = 9.858689999999999
$ws.Range("H17").Value = 19.71738
$ws.Range("I17").Value = 0.04562679349910256
$ws.Range("J17").Value = 0.03088763017394275
$ws.Range("M17").Value = 162.7225033333333
$ws.Range("N17").Value = 488.16751
$ws.Range("O17").Value = 0.5231437953541009
$ws.Range("P17").Value = 0.5247717033381212
$ws.Range("Q17").Value = 1604.2307163873
$ws.Range("R17").Value = 9625.384298323799
$ws.Range("S17").Value = 0.02386937392095833
$ws.Range("T17").Value = 0.01620895429845788
$ws.Range("G18").Value = 9.858689999999999
$ws.Range("H18").Value = 19.71738
$ws.Range("I18").Value = 0.04562679349910256
$ws.Range("J18").Value = 0.03088763017394275
$ws.Range("O18").Value = 0.0009322191998643353
$ws.Range("P18").Value = 0.0009351200601857102
$ws.Range("Q18").Value = 2.85866847339
$ws.Range("R18").Value = 17.15201084034
$ws.Range("S18").Value = 0.00004253417292810865
$ws.Range("T18").Value = 0.0000288836425872513
$ws.Range("G19").Value = 9.858689999999999
$ws.Range("H19").Value = 19.71738
$ws.Range("I19").Value = 0.04562679349910256
$ws.Range("J19").Value = 0.03088763017394275
$ws.Range("M19").Value = 61.580654
$ws.Range("N19").Value = 184.741962
$ws.Range("O19").Value = 0.1979783766474813
$ws.Range("P19").Value = 0.1985944416431287
$ws.Range("Q19").Value = 607.10457778326
$ws.Range("R19").Value = 3642.62746669956
$ws.Range("S19").Value = 0.00903311850858218
$ws.Range("T19").Value = 0.006134111668073615
$ws.Range("G20").Value = 9.858689999999999
$ws.Range("H20").Value = 19.71738
$ws.Range("I20").Value = 0.04562679349910256
$ws.Range("J20").Value = 0.03088763017394275
$ws.Range("M20").Value = 2.8947245
$ws.Range("N20").Value = 5.789449
$ws.Range("O20").Value = 0.009306378223129816
$ws.Range("P20").Value = 0.00622355841157717
$ws.Range("Q20").Value = 28.538191480905
$ws.Range("R20").Value = 114.15276592362
$ws.Range("S20").Value = 0.0004246201974112891
$ws.Range("T20").Value = 0.0001922309705827262
$ws.Range("G21").Value = 9.858689999999999
$ws.Range("H21").Value = 19.71738
$ws.Range("I21").Value = 0.04562679349910256
$ws.Range("J21").Value = 0.03088763017394275
$ws.Range("M21").Value = 83.559527
$ws.Range("N21").Value = 250.678581
$ws.Range("O21").Value = 0.2686392305754237
$ws.Range("P21").Value = 0.2694751765469873
$ws.Range("Q21").Value = 823.78747323963
$ws.Range("R21").Value = 4942.72483943778
$ws.Range("S21").Value = 0.01225714669922266
$ws.Range("T21").Value = 0.008323449594241274
$ws.Range("G22").Value = 1.523862
$ws.Range("H22").Value = 4.571586
$ws.Range("I22").Value = 0.007052553310341378
$ws.Range("J22").Value = 0.007161471639557297
$ws.Range("M22").Value = 162.7225033333333
$ws.Range("N22").Value = 488.16751
$ws.Range("O22").Value = 0.5231437953541009
$ws.Range("P22").Value = 0.5247717033381212
$ws.Range("Q22").Value = 247.96663937454
$ws.Range("R22").Value = 2231.69975437086
$ws.Range("S22").Value = 0.003689499505709117
$ws.Range("T22").Value = 0.00375813767069813
$ws.Range("G23").Value = 1.523862
$ws.Range("H23").Value = 4.571586
$ws.Range("I23").Value = 0.007052553310341378
$ws.Range("J23").Value = 0.007161471639557297
$ws.Range("O23").Value = 0.0009322191998643353
$ws.Range("P23").Value = 0.0009351200601857102
$ws.Range("Q23").Value = 0.441865628922
$ws.Range("R23").Value = 3.976790660298
$ws.Range("S23").Value = 0.000006574525603967009
$ws.Range("T23").Value = 0.000006696835790601076
$ws.Range("G24").Value = 1.523862
$ws.Range("H24").Value = 4.571586
$ws.Range("I24").Value = 0.007052553310341378
$ws.Range("J24").Value = 0.007161471639557297
$ws.Range("M24").Value = 61.580654
$ws.Range("N24").Value = 184.741962
$ws.Range("O24").Value = 0.1979783766474813
$ws.Range("P24").Value = 0.1985944416431287
$ws.Range("Q24").Value = 93.84041856574801
$ws.Range("R24").Value = 844.563767091732
$ws.Range("S24").Value = 0.001396253055601206
$ws.Range("T24").Value = 0.001422228461600983
$ws.Range("G25").Value = 1.523862
$ws.Range("H25").Value = 4.571586
$ws.Range("I25").Value = 0.007052553310341378
$ws.Range("J25").Value = 0.007161471639557297
$ws.Range("M25").Value = 2.8947245
$ws.Range("N25").Value = 5.789449
$ws.Range("O25").Value = 0.009306378223129816
$ws.Range("P25").Value = 0.00622355841157717
$ws.Range("Q25").Value = 4.411160666019001
$ws.Range("R25").Value = 26.466963996114
$ws.Range("S25").Value = 0.00006563372854482309
$ws.Range("T25").Value = 0.00004456983706163816
$ws.Range("G26").Value = 1.523862
$ws.Range("H26").Value = 4.571586
$ws.Range("I26").Value = 0.007052553310341378
$ws.Range("J26").Value = 0.007161471639557297
$ws.Range("M26").Value = 83.559527
$ws.Range("N26").Value = 250.678581
$ws.Range("O26").Value = 0.2686392305754237
$ws.Range("P26").Value = 0.2694751765469873
$ws.Range("Q26").Value = 127.333187933274
$ws.Range("R26").Value = 1145.998691399466
$ws.Range("S26").Value = 0.001894592494882265
$ws.Range("T26").Value = 0.001929838834405945
